$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$sub3 = [char]0x2083

Set-CellText "D2" "67.633.68"
Set-CellText "E2" "  +3.40%  "
Set-CellText "D3" "3.769.05"
Set-CellText "E3" "  +7.76%  "
Set-CellText "E4" "  +0.08%  "
Set-CellText "D5" "418.96"
Set-CellText "E5" "  +0.40%  "
Set-CellText "D6" "132.25"
Set-CellText "E6" "  +1.06%  "
Set-CellText "D7" "3.751.26"
Set-CellText "E7" "  +7.55%  "
Set-CellText "D8" "0.645"
Set-CellText "E8" "  -1.68%  "
Set-CellText "E9" "  +0.04%  "
Set-CellText "D10" "0.765"
Set-CellText "E10" "  -2.00%  "
Set-CellText "D11" "0.184"
Set-CellText "E11" "  +13.09%  "
Set-CellText "D12" "0.0000397"
Set-CellText "E12" "  +48.51%  "
Set-CellText "D13" "42.37"
Set-CellText "E13" "  -2.08%  "
Set-CellText "D14" "10.48"
Set-CellText "E14" "  +4.80%  "
Set-CellText "D15" "4.373.23"
Set-CellText "E15" "  +7.96%  "
Set-CellText "E16" "  -0.60%  "
Set-CellText "D17" "3.767.62"
Set-CellText "E17" "  +7.78%  "
Set-CellText "D18" "20.38"
Set-CellText "E18" "  -0.22%  "
Set-CellText "D19" "13.22"
Set-CellText "E19" "  +2.59%  "
Set-CellText "E20" "  +2.09%  "
Set-CellText "D21" "67.524.65"
Set-CellText "E21" "  +3.76%  "
Set-CellText "D22" "441.21"
Set-CellText "E22" "  -1.01%  "
Set-CellText "D23" "15.58"
Set-CellText "E23" "  +18.52%  "
Set-CellText "D24" "89.89"
Set-CellText "E24" "  +0.33%  "
Set-CellText "D25" "3.08"
Set-CellText "E25" "  -5.05%  "
Set-CellText "D26" "38.16"
Set-CellText "E26" "  +12.34%  "
Set-CellText "D27" "3.31"
Set-CellText "E27" "  -1.57%  "
Set-CellText "D28" "10.03"
Set-CellText "E28" "  +1.40%  "
Set-CellText "D29" "5.08"
Set-CellText "E29" "  +5.20%  "
Set-CellText "E30" "  +5.79%  "
Set-CellText "E31" "  +0.65%  "
Set-CellText "D32" "2.77"
Set-CellText "E32" "  +1.12%  "
Set-CellText "D33" "7.13"
Set-CellText "E33" "  -3.45%  "
Set-CellText "D34" "0.163"
Set-CellText "E34" "  +0.09%  "
Set-CellText "D35" "41.04"
Set-CellText "E35" "  +4.13%  "
Set-CellText "D36" "57.93"
Set-CellText "E36" "  +0.89%  "
Set-CellText "D37" "0.998"
Set-CellText "E37" "  +0.05%  "
Set-CellText "D38" "0.0486"
Set-CellText "E38" "  -3.44%  "
Set-CellText "D39" "2.99"
Set-CellText "E39" "  +28.74%  "
Set-CellText "B40" "PEPE"
Set-CellText "C40" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText "D40" ("0.0$sub3" + "0708")
Set-CellText "E40" "  -2.51%  "
Set-CellText "B41" "Stellar"
Set-CellText "C41" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText "D41" "0.147"
Set-CellText "E41" "  +0.54%  "
Set-CellText "D42" "0.996"
Set-CellText "E42" "  -0.02%  "
Set-CellText "D43" "27.92"
Set-CellText "E43" "  +29.79%  "
Set-CellText "E44" "  +3.98%  "
Set-CellText "D45" "147.88"
Set-CellText "E45" "  +0.39%  "
Set-CellText "E46" "  +24.62%  "
Set-CellText "D47" "2.09"
Set-CellText "E47" "  +4.76%  "
Set-CellText "D48" "2.87"
Set-CellText "E48" "  -4.50%  "
Set-CellText "B49" "NEARProtocol"
Set-CellText "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-CellText "D49" "4.28"
Set-CellText "E49" "  -4.83%  "
Set-CellText "B50" "WEMIXToken"
Set-CellText "C50" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText "D50" "2.59"
Set-CellText "E50" "  -7.22%  "
Set-CellText "D51" "0.302"
Set-CellText "E51" "  -2.70%  "
